# The deck ships two theme parts:
#   - ppt/theme/theme1.xml  -> used by the Slide Master (and therefore every
#     slide). It currently carries the "Integral" / "Red Violet" palette and
#     must be repainted with the standard "Office" palette.
#   - ppt/theme/theme2.xml  -> used only by the Notes Master. It currently
#     carries the "Office" palette and (per the target diff) would become
#     the "Integral" / "Red Violet" palette.
#
# PowerPoint's object model has no bulk "swap these two themes" call, so each
# of the twelve theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink -- PpThemeColorSchemeIndex order 1..12) is repainted individually
# through ThemeColorScheme.Item(n).RGB, which is the supported read/write
# surface for theme colours.

function HexToComRGB($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

function SetThemeColors($theme, $hexColors) {
    $tcs = $theme.ThemeColorScheme
    for ($i = 1; $i -le $hexColors.Count; $i++) {
        $tcs.Item($i).RGB = HexToComRGB $hexColors[$i - 1]
    }
}

$p = $ppt.ActivePresentation

# PpThemeColorSchemeIndex order: dk1, lt1, dk2, lt2,
# accent1, accent2, accent3, accent4, accent5, accent6, hlink, folHlink
$officeColors = @(
    "000000", "FFFFFF", "44546A", "E7E6E6",
    "5B9BD5", "ED7D31", "A5A5A5", "FFC000",
    "4472C4", "70AD47", "0563C1", "954F72"
)

# Slide master theme (theme1.xml): Integral -> Office
SetThemeColors $p.SlideMaster.Theme $officeColors
